$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set values in the exact order needed so new shared-strings are created
# in the same sequence as the target workbook (37..41).
$ws.Range("K2").Value = "None"
$ws.Range("K11").Value = "Come Back To"
$ws.Range("K4").Value = "Done, most weak signal"
$ws.Range("J6").Value = "weak+D6:J6, hiss"
$ws.Range("K19").Value = "None, Looks like burst will be populated"

# Remaining "Done" / "None" cells (reuse existing shared strings)
$ws.Range("K5").Value = "Done"
$ws.Range("K6").Value = "Done"
$ws.Range("K7").Value = "Done"
$ws.Range("K8").Value = "Done"
$ws.Range("K9").Value = "None"
$ws.Range("K10").Value = "None"
$ws.Range("K12").Value = "Done"
$ws.Range("K14").Value = "Done"
$ws.Range("K15").Value = "None"
$ws.Range("K16").Value = "Done"
$ws.Range("K17").Value = "Done"
$ws.Range("K18").Value = "Done"
$ws.Range("K20").Value = "Done"
